# Add two new columns, I (I0) and J (IF), to the existing data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, bordered, centered) used by the other
# header cells (B1:H1) by copying the style from H1 onto the new headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-17) -------------------------------------------------
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(2, 10).Value = 4

$ws.Cells.Item(3, 9).Value = 8
$ws.Cells.Item(3, 10).Value = 9

$ws.Cells.Item(4, 9).Value = 7
$ws.Cells.Item(4, 10).Value = 9

$ws.Cells.Item(5, 9).Value = 3
$ws.Cells.Item(5, 10).Value = 6

$ws.Cells.Item(6, 9).Value = 3
$ws.Cells.Item(6, 10).Value = 6

$ws.Cells.Item(7, 9).Value = 6
$ws.Cells.Item(7, 10).Value = 8

$ws.Cells.Item(8, 9).Value = 2
$ws.Cells.Item(8, 10).Value = 6

$ws.Cells.Item(9, 9).Value = 3
$ws.Cells.Item(9, 10).Value = 7

$ws.Cells.Item(10, 9).Value = 3
$ws.Cells.Item(10, 10).Value = 6

$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(11, 10).Value = 8

$ws.Cells.Item(12, 9).Value = 1
$ws.Cells.Item(12, 10).Value = 4

$ws.Cells.Item(13, 9).Value = 5
$ws.Cells.Item(13, 10).Value = 7

$ws.Cells.Item(14, 9).Value = 6
$ws.Cells.Item(14, 10).Value = 7

$ws.Cells.Item(15, 9).Value = 3
$ws.Cells.Item(15, 10).Value = 3

$ws.Cells.Item(16, 9).Value = 7
$ws.Cells.Item(16, 10).Value = 7

$ws.Cells.Item(17, 9).Value = 1
$ws.Cells.Item(17, 10).Value = 1
